# edit.ps1 - applies "New crime data collected" weekly update
# to the NYPD CompStat report workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header banner text: Volume/Number and the reporting week's date range.
#    These are rich-text cells; replace only the digits/date substrings so
#    the surrounding label text is untouched.
# ---------------------------------------------------------------------------

# A8: "Volume 29   Number  44" -> "Volume 29   Number  45"
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "45"

# C9: "Report Covering the Week  10/31/2022  Through  11/6/2022"
#  -> "Report Covering the Week  11/7/2022  Through  11/13/2022"
$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "11/7/2022"
$c9.Characters(47, 9).Text = "11/13/2022"

# ---------------------------------------------------------------------------
# 2. Cells that flip from a text placeholder ("0" / "***.*") to a real
#    number now that there is enough data to report an actual figure.
#    Re-apply the numeric format used by sibling cells in the same table
#    before writing the value, so the style matches the rest of the column.
# ---------------------------------------------------------------------------

$ws.Range("C14").NumberFormat = $ws.Range("D14").NumberFormat
$ws.Range("C14").Value = 1

$ws.Range("D28").NumberFormat = $ws.Range("C28").NumberFormat
$ws.Range("D28").Value = 1

$ws.Range("E28").NumberFormat = $ws.Range("H28").NumberFormat
$ws.Range("E28").Value = 300

$ws.Range("D29").NumberFormat = $ws.Range("C29").NumberFormat
$ws.Range("D29").Value = 1

$ws.Range("E29").NumberFormat = $ws.Range("H29").NumberFormat
$ws.Range("E29").Value = 100

# ---------------------------------------------------------------------------
# 3. Refreshed weekly crime-complaint figures (rows 14-30): new week-to-date
#    counts, 28-day counts, year-to-date counts and the derived percentages.
# ---------------------------------------------------------------------------

$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = 2
$ws.Cells.Item(14,7).Value = 2
$ws.Cells.Item(14,8).Value = 0
$ws.Cells.Item(14,9).Value = 26
$ws.Cells.Item(14,10).Value = 23
$ws.Cells.Item(14,11).Value = 13.043478260869
$ws.Cells.Item(14,12).Value = 52.941176470588
$ws.Cells.Item(14,13).Value = 73.333333333333
$ws.Cells.Item(14,14).Value = -64.864864864864
$ws.Cells.Item(15,3).Value = 3
$ws.Cells.Item(15,4).Value = 1
$ws.Cells.Item(15,6).Value = 14
$ws.Cells.Item(15,7).Value = 10
$ws.Cells.Item(15,8).Value = 40
$ws.Cells.Item(15,9).Value = 157
$ws.Cells.Item(15,10).Value = 136
$ws.Cells.Item(15,11).Value = 15.441176470588
$ws.Cells.Item(15,12).Value = 30.833333333333
$ws.Cells.Item(15,13).Value = 63.541666666666
$ws.Cells.Item(15,14).Value = -18.652849740932
$ws.Cells.Item(16,3).Value = 55
$ws.Cells.Item(16,4).Value = 46
$ws.Cells.Item(16,5).Value = 19.565217391304
$ws.Cells.Item(16,6).Value = 180
$ws.Cells.Item(16,7).Value = 186
$ws.Cells.Item(16,8).Value = -3.225806451612
$ws.Cells.Item(16,9).Value = 1877
$ws.Cells.Item(16,10).Value = 1442
$ws.Cells.Item(16,11).Value = 30.166435506241
$ws.Cells.Item(16,12).Value = 77.914691943128
$ws.Cells.Item(16,13).Value = 61.67097329888
$ws.Cells.Item(16,14).Value = -79.650910667823
$ws.Cells.Item(17,3).Value = 39
$ws.Cells.Item(17,4).Value = 39
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = 176
$ws.Cells.Item(17,7).Value = 170
$ws.Cells.Item(17,8).Value = 3.529411764705
$ws.Cells.Item(17,9).Value = 1846
$ws.Cells.Item(17,10).Value = 1677
$ws.Cells.Item(17,11).Value = 10.077519379845
$ws.Cells.Item(17,12).Value = 46.275752773375
$ws.Cells.Item(17,13).Value = 70.609981515711
$ws.Cells.Item(17,14).Value = -39.455559199737
$ws.Cells.Item(18,3).Value = 49
$ws.Cells.Item(18,4).Value = 41
$ws.Cells.Item(18,5).Value = 19.512195121951
$ws.Cells.Item(18,6).Value = 199
$ws.Cells.Item(18,7).Value = 207
$ws.Cells.Item(18,8).Value = -3.864734299516
$ws.Cells.Item(18,9).Value = 2565
$ws.Cells.Item(18,10).Value = 1795
$ws.Cells.Item(18,11).Value = 42.896935933147
$ws.Cells.Item(18,12).Value = 13.395225464191
$ws.Cells.Item(18,13).Value = 65.164198325821
$ws.Cells.Item(18,14).Value = -75.659517935092
$ws.Cells.Item(19,3).Value = 205
$ws.Cells.Item(19,4).Value = 255
$ws.Cells.Item(19,5).Value = -19.607843137254
$ws.Cells.Item(19,6).Value = 1035
$ws.Cells.Item(19,7).Value = 898
$ws.Cells.Item(19,8).Value = 15.256124721603
$ws.Cells.Item(19,9).Value = 10189
$ws.Cells.Item(19,10).Value = 6793
$ws.Cells.Item(19,11).Value = 49.992639481819
$ws.Cells.Item(19,12).Value = 66.052803129074
$ws.Cells.Item(19,13).Value = 10.786125910623
$ws.Cells.Item(19,14).Value = -64.644852354349
$ws.Cells.Item(20,3).Value = 13
$ws.Cells.Item(20,5).Value = 18.181818181818
$ws.Cells.Item(20,6).Value = 54
$ws.Cells.Item(20,7).Value = 51
$ws.Cells.Item(20,8).Value = 5.882352941176
$ws.Cells.Item(20,9).Value = 606
$ws.Cells.Item(20,10).Value = 532
$ws.Cells.Item(20,11).Value = 13.90977443609
$ws.Cells.Item(20,12).Value = 33.186813186813
$ws.Cells.Item(20,13).Value = 72.15909090909
$ws.Cells.Item(20,14).Value = -88.931506849315
$ws.Cells.Item(21,3).Value = 365
$ws.Cells.Item(21,4).Value = 394
$ws.Cells.Item(21,5).Value = -7.36040609137
$ws.Cells.Item(21,6).Value = 1660
$ws.Cells.Item(21,7).Value = 1524
$ws.Cells.Item(21,8).Value = 8.923884514435
$ws.Cells.Item(21,9).Value = 17266
$ws.Cells.Item(21,10).Value = 12398
$ws.Cells.Item(21,11).Value = 39.264397483465
$ws.Cells.Item(21,12).Value = 52.701866100645
$ws.Cells.Item(21,13).Value = 28.314506539833
$ws.Cells.Item(21,14).Value = -69.905180227288
$ws.Cells.Item(22,3).Value = 8
$ws.Cells.Item(22,4).Value = 16
$ws.Cells.Item(22,5).Value = -50
$ws.Cells.Item(22,6).Value = 58
$ws.Cells.Item(22,7).Value = 51
$ws.Cells.Item(22,8).Value = 13.725490196078
$ws.Cells.Item(22,9).Value = 583
$ws.Cells.Item(22,10).Value = 429
$ws.Cells.Item(22,11).Value = 35.897435897435
$ws.Cells.Item(22,12).Value = 35.266821345707
$ws.Cells.Item(22,13).Value = 23.255813953488
$ws.Cells.Item(23,3).Value = 6
$ws.Cells.Item(23,4).Value = 4
$ws.Cells.Item(23,5).Value = 50
$ws.Cells.Item(23,6).Value = 23
$ws.Cells.Item(23,7).Value = 33
$ws.Cells.Item(23,8).Value = -30.30303030303
$ws.Cells.Item(23,9).Value = 381
$ws.Cells.Item(23,10).Value = 433
$ws.Cells.Item(23,11).Value = -12.009237875288
$ws.Cells.Item(23,12).Value = -2.557544757033
$ws.Cells.Item(23,13).Value = 22.508038585209
$ws.Cells.Item(24,3).Value = 405
$ws.Cells.Item(24,4).Value = 368
$ws.Cells.Item(24,5).Value = 10.054347826087
$ws.Cells.Item(24,6).Value = 1713
$ws.Cells.Item(24,7).Value = 1504
$ws.Cells.Item(24,8).Value = 13.896276595744
$ws.Cells.Item(24,9).Value = 19058
$ws.Cells.Item(24,10).Value = 12681
$ws.Cells.Item(24,11).Value = 50.28783218989
$ws.Cells.Item(24,12).Value = 75.698349774131
$ws.Cells.Item(24,13).Value = 28.979426096372
$ws.Cells.Item(25,3).Value = 101
$ws.Cells.Item(25,4).Value = 84
$ws.Cells.Item(25,5).Value = 20.238095238095
$ws.Cells.Item(25,6).Value = 419
$ws.Cells.Item(25,7).Value = 373
$ws.Cells.Item(25,8).Value = 12.332439678284
$ws.Cells.Item(25,9).Value = 4154
$ws.Cells.Item(25,10).Value = 3479
$ws.Cells.Item(25,11).Value = 19.402127048002
$ws.Cells.Item(25,12).Value = 52.664461595001
$ws.Cells.Item(25,13).Value = 31.082360366046
$ws.Cells.Item(26,3).Value = 6
$ws.Cells.Item(26,5).Value = 50
$ws.Cells.Item(26,6).Value = 22
$ws.Cells.Item(26,8).Value = 22.222222222222
$ws.Cells.Item(26,9).Value = 245
$ws.Cells.Item(26,10).Value = 206
$ws.Cells.Item(26,11).Value = 18.932038834951
$ws.Cells.Item(26,12).Value = 29.629629629629
$ws.Cells.Item(27,3).Value = 13
$ws.Cells.Item(27,4).Value = 17
$ws.Cells.Item(27,5).Value = -23.529411764705
$ws.Cells.Item(27,6).Value = 74
$ws.Cells.Item(27,7).Value = 57
$ws.Cells.Item(27,8).Value = 29.824561403508
$ws.Cells.Item(27,9).Value = 861
$ws.Cells.Item(27,10).Value = 756
$ws.Cells.Item(27,11).Value = 13.888888888888
$ws.Cells.Item(27,12).Value = 58.856088560885
$ws.Cells.Item(28,3).Value = 4
$ws.Cells.Item(28,6).Value = 8
$ws.Cells.Item(28,7).Value = 2
$ws.Cells.Item(28,8).Value = 300
$ws.Cells.Item(28,9).Value = 49
$ws.Cells.Item(28,10).Value = 37
$ws.Cells.Item(28,11).Value = 32.432432432432
$ws.Cells.Item(28,12).Value = 53.125
$ws.Cells.Item(28,13).Value = 36.111111111111
$ws.Cells.Item(28,14).Value = -58.474576271186
$ws.Cells.Item(29,3).Value = 2
$ws.Cells.Item(29,6).Value = 5
$ws.Cells.Item(29,8).Value = 150
$ws.Cells.Item(29,9).Value = 41
$ws.Cells.Item(29,10).Value = 33
$ws.Cells.Item(29,11).Value = 24.242424242424
$ws.Cells.Item(29,12).Value = 46.428571428571
$ws.Cells.Item(29,13).Value = 57.692307692307
$ws.Cells.Item(29,14).Value = -61.320754716981
$ws.Cells.Item(30,3).Value = 1
$ws.Cells.Item(30,4).Value = 1
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = 10
$ws.Cells.Item(30,8).Value = 66.666666666666
$ws.Cells.Item(30,9).Value = 138
$ws.Cells.Item(30,10).Value = 129
$ws.Cells.Item(30,11).Value = 6.976744186046
$ws.Cells.Item(30,12).Value = 187.5

